$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.106.18"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.638.00"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'216.60"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'19.91"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.869.52"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "1.639.78"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "'66.77"
$ws.Range("D17").Value = "27.138.53"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "'216.78"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'6.94"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'2.53"
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "'146.78"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'7.41"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'15.66"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").Value = "'0.0508"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").Value = "1.306.67"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").Value = "'0.0175"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").Value = "'0.544"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'0.811"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +5.73%  "
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "1.779.01"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'61.67"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "'91.44"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").Value = "'0.0511"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "'7.65"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").Value = "'0.0961"
$ws.Range("E51").Value = "  -0.18%  "
